$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "57.834.54"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "533.17"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "139.00"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.126.71"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("D13").Value = "3.666.49"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "25.90"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "57.930.61"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "3.120.61"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("E20").Value = "  +3.13%  "
$ws.Range("D21").Value = "8.08"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "367.65"
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "5.66"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "69.24"
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("D26").Value = "0.506"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").Value = "21.43"
$ws.Range("E33").Value = "  +3.99%  "
$ws.Range("D34").Value = "5.17"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").Value = "159.47"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  +5.89%  "
$ws.Range("D39").Value = "25.47"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  +5.02%  "
$ws.Range("D41").Value = "0.0671"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "2.521.10"
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("D43").Value = "4.09"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "37.83"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("D50").Value = "19.74"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  -0.57%  "
